# Update "total_risk" (column R) and "total_risk_resp" (column S) values
# on Sheet 1 for rows 2-15, per the transition rule 5 and 10 mi radius
# script/output updates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$updates = @(
    @{ Row = 2;  R = 35;    S = 0.45 },
    @{ Row = 3;  R = 50;    S = 0.6 },
    @{ Row = 4;  R = 40;    S = 0.442857142857143 },
    @{ Row = 5;  R = 30;    S = 0.485714285714286 },
    @{ Row = 6;  R = 30;    S = 0.3 },
    @{ Row = 7;  R = 30;    S = 0.4 },
    @{ Row = 8;  R = 20;    S = 0.2 },
    @{ Row = 9;  R = 90;    S = 0.5 },
    @{ Row = 10; R = 107.5; S = 0.525 },
    @{ Row = 11; R = 40;    S = 0.5 },
    @{ Row = 12; R = 30;    S = 0.35 },
    @{ Row = 13; R = 40;    S = 0.4 },
    @{ Row = 14; R = 20;    S = 0.2 },
    @{ Row = 15; R = 30;    S = 0.4 }
)

foreach ($u in $updates) {
    $ws.Range("R" + $u.Row).Value = $u.R
    $ws.Range("S" + $u.Row).Value = $u.S
}
